$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tidsregistrering")

# ---------------------------------------------------------------------------
# 1) Make room for the two brand-new activity rows on 2017-03-07 (row 10 already
#    exists as a placeholder date row; rows 11 & 12 need to be inserted so the
#    old rows 11-27 become 13-29).
# ---------------------------------------------------------------------------
$ws.Rows.Item(11).Insert()
$ws.Rows.Item(11).Insert()

# ---------------------------------------------------------------------------
# 2) Participant column (C) - "Nada H. A. Omer" is filled in for every
#    activity row (4 through 12). The remaining date/placeholder rows down to
#    28 just pick up the bold "Navn:"-style formatting (matching the extended
#    Deltagere validation range) without an actual value.
# ---------------------------------------------------------------------------
$ws.Range("C4:C12").Value = "Nada H. A. Omer"
$ws.Range("C13:C28").Font.Bold = $true

# ---------------------------------------------------------------------------
# 3) Role column moves from D to E for the existing 2017-03-06 activities, and
#    the new 2017-03-07 activities get their role filled directly into E.
# ---------------------------------------------------------------------------
$ws.Range("D5:D8").ClearContents()
$ws.Range("E4:E9").Value = "System Analyst"
$ws.Range("E10:E11").Value = "Requirements Specifier"
$ws.Range("E12").Value = "Project Manager"

# ---------------------------------------------------------------------------
# 4) Fill in the three new activity rows for 2017-03-07.
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = "NO"
$ws.Range("F10").Value = "Visionsdokument + FURPS"
$ws.Range("G10").Value = 0.36458333333333331
$ws.Range("H10").Value = 0.4375
$ws.Range("I10").Value = "1 time : 45 min."

$ws.Range("A11").Value = $ws.Range("A10").Value
$ws.Range("B11").Value = "NO"
$ws.Range("F11").Value = "Formel Use case 5 + Domæne model for UC 5"
$ws.Range("G11").Value = 0.4375
$ws.Range("H11").Value = 0.47916666666666669
$ws.Range("I11").Value = "1 time : 00"

$ws.Range("A12").Value = $ws.Range("A10").Value
$ws.Range("B12").Value = "NO"
$ws.Range("F12").Value = "Iterationsplan"
$ws.Range("G12").Value = 0.52083333333333337
$ws.Range("H12").Value = 0.63888888888888895
$ws.Range("I12").Value = "2 time : 50"

# ---------------------------------------------------------------------------
# 5) Data validation ranges grow to cover the newly used rows.
# ---------------------------------------------------------------------------
$ws.Range("E3:E107").Validation.Delete()
$ws.Range("E3:E109").Validation.Add(3, 1, 1, "=GyldigeRoller")

$ws.Range("C3").Validation.Delete()
$ws.Range("C3:C28").Validation.Add(3, 1, 1, "=Deltagere")

# ---------------------------------------------------------------------------
# 6) Column width tweaks to better fit the new content.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 6.42578125
$ws.Columns.Item(3).ColumnWidth = 17.85546875
$ws.Columns.Item(4).ColumnWidth = 7.140625
$ws.Columns.Item(5).ColumnWidth = 22.28515625

# ---------------------------------------------------------------------------
# 7) Keep the author's last-used selection.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("I15").Select()
